$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RO & CO Hearing Allocation")
$ws.Rows.Item(4).Delete()
